$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-21 23:48:38"
$ws.Range("E3").Value = "2026-02-21 23:48:40"
$ws.Range("E4").Value = "2026-02-21 23:48:43"
$ws.Range("H4").Value = "75%"
$ws.Range("O4").Value = "8.9 °C"
$ws.Range("E5").Value = "2026-02-21 23:48:46"
$ws.Range("H5").Value = "37%"
$ws.Range("E6").Value = "2026-02-21 23:48:48"
$ws.Range("O6").Value = "10.5 °C"
$ws.Range("E7").Value = "2026-02-21 23:48:51"
$ws.Range("E8").Value = "2026-02-21 23:48:54"
$ws.Range("E9").Value = "2026-02-21 23:48:56"
$ws.Range("H9").Value = "61%"
$ws.Range("O9").Value = "12.6 °C"
$ws.Range("E10").Value = "2026-02-21 23:48:59"
$ws.Range("O10").Value = "8.1 °C"
$ws.Range("E11").Value = "2026-02-21 23:49:01"
$ws.Range("H11").Value = "56%"
$ws.Range("O11").Value = "8.3 °C"
$ws.Range("E12").Value = "2026-02-21 23:49:03"
$ws.Range("N12").Value = "6.1 °C 23:10 TU"
$ws.Range("O12").Value = "12.1 °C"
$ws.Range("E13").Value = "2026-02-21 23:49:06"
$ws.Range("H13").Value = "62%"
$ws.Range("O13").Value = "5.0 °C"
$ws.Range("E14").Value = "2026-02-21 23:49:08"
$ws.Range("H14").Value = "73%"
$ws.Range("O14").Value = "10.9 °C"
$ws.Range("E15").Value = "2026-02-21 23:49:11"
$ws.Range("O15").Value = "12.3 °C"
$ws.Range("E16").Value = "2026-02-21 23:49:13"
$ws.Range("E17").Value = "2026-02-21 23:49:16"
$ws.Range("E18").Value = "2026-02-21 23:49:19"
$ws.Range("H18").Value = "77%"
$ws.Range("O18").Value = "8.2 °C"
$ws.Range("E19").Value = "2026-02-21 23:49:21"
$ws.Range("E20").Value = "2026-02-21 23:49:24"
$ws.Range("E21").Value = "2026-02-21 23:49:26"
$ws.Range("O21").Value = "7.4 °C"
$ws.Range("E22").Value = "2026-02-21 23:49:29"
$ws.Range("O22").Value = "2.1 °C"
$ws.Range("E23").Value = "2026-02-21 23:49:31"
$ws.Range("E24").Value = "2026-02-21 23:49:34"
$ws.Range("H24").Value = "85%"
$ws.Range("E25").Value = "2026-02-21 23:49:37"
$ws.Range("E26").Value = "2026-02-21 23:49:39"
$ws.Range("E27").Value = "2026-02-21 23:49:42"
$ws.Range("E28").Value = "2026-02-21 23:49:45"
$ws.Range("O28").Value = "7.9 °C"
$ws.Range("E29").Value = "2026-02-21 23:49:47"
$ws.Range("H29").Value = "70%"
$ws.Range("E30").Value = "2026-02-21 23:49:50"
$ws.Range("H30").Value = "70%"
$ws.Range("O30").Value = "11.1 °C"
$ws.Range("E31").Value = "2026-02-21 23:49:53"
$ws.Range("E32").Value = "2026-02-21 23:49:56"
$ws.Range("O32").Value = "4.3 °C"
$ws.Range("E33").Value = "2026-02-21 23:49:58"
$ws.Range("J33").Value = "1030.7 hPa"
$ws.Range("O33").Value = "6.4 °C"
$ws.Range("E34").Value = "2026-02-21 23:50:01"
$ws.Range("O34").Value = "4.1 °C"
$ws.Range("E35").Value = "2026-02-21 23:50:03"
$ws.Range("E36").Value = "2026-02-21 23:50:06"
$ws.Range("H36").Value = "62%"
$ws.Range("N36").Value = "5.9 °C 23:20 TU"
$ws.Range("O36").Value = "12.7 °C"
$ws.Range("E37").Value = "2026-02-21 23:50:09"
$ws.Range("H37").Value = "76%"
$ws.Range("J37").Value = "1031.9 hPa"
$ws.Range("O37").Value = "5.3 °C"
$ws.Range("E38").Value = "2026-02-21 23:50:12"
$ws.Range("K38").Value = "15.5 MJ/m2"
$ws.Range("E39").Value = "2026-02-21 23:50:14"
$ws.Range("O39").Value = "2.5 °C"
$ws.Range("E40").Value = "2026-02-21 23:50:17"
$ws.Range("H40").Value = "56%"
$ws.Range("J40").Value = "1031.0 hPa"
$ws.Range("O40").Value = "8.1 °C"
$ws.Range("E41").Value = "2026-02-21 23:50:19"
$ws.Range("E42").Value = "2026-02-21 23:50:22"
$ws.Range("O42").Value = "10.4 °C"
$ws.Range("E43").Value = "2026-02-21 23:50:25"
$ws.Range("O43").Value = "7.0 °C"
$ws.Range("E44").Value = "2026-02-21 23:50:28"
$ws.Range("E45").Value = "2026-02-21 23:50:31"
$ws.Range("K45").Value = "13.6 MJ/m2"
$ws.Range("E46").Value = "2026-02-21 23:50:33"
$ws.Range("O46").Value = "9.2 °C"
